# Updates VTQaZ sheet: allow hydrogen FCVs to qualify as zero emission.
# Row 8 ("hydrogen vehicle") values in columns B:AF change from 0 to 1.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VTQaZ")

# Set the row of values (B8:AF8) to 1 (was 0)
$ws.Range("B8:AF8").Value = 1

# Select the range and activate it to match the saved selection state
$ws.Activate()
$ws.Range("B8:AF8").Select()
